$d = $word.ActiveDocument

$old = "年キャンペーン期間 対象：うしかい座星座 2022: 5月14日〜23日、6月13日〜22日、7月12日〜21日"
$new = " ：2022年キャンペーン期間 (対象：うしかい座星座)：、5月14日〜23日、6月13日〜22日、7月12日〜21日"

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
